# A sale of 10 units of "Amul Butter (500g)" was recorded:
#   - Sales sheet gets a new row (row 3) describing the sale.
#   - Products sheet stock for that item drops from 29 to 19 (29 - 10)
#     and its last_updated timestamp is refreshed to match the sale time.

$wb = $excel.ActiveWorkbook

$products = $wb.Worksheets.Item("Products")
$sales    = $wb.Worksheets.Item("Sales")

# --- Products sheet: decrement stock for the sold item, refresh timestamp ---
$products.Range("F2").Value = 19
$products.Range("L2").Value = "2025-09-23T13:02:08.605Z"

# --- Sales sheet: append the new sale as row 3 ---
$sales.Range("A3").Value = "07291cea-a90f-4352-b915-1121e63dcb59"
$sales.Range("B3").Value = 10

# Columns that hold numeric-looking values but are stored as text in this
# workbook (matching the sibling cells already on row 2) need NumberFormat
# forced to text first, otherwise they'd be written back as real numbers.
$sales.Range("C3").NumberFormat = "@"
$sales.Range("C3").Value = "500"

$sales.Range("D3").NumberFormat = "@"
$sales.Range("D3").Value = ""

$sales.Range("E3").Value = "Cash"
$sales.Range("F3").Value = "Admin"

$sales.Range("G3").NumberFormat = "@"
$sales.Range("G3").Value = ""

$sales.Range("H3").Value = "e66c2199-14fc-46af-b2ee-ed7812584635"
$sales.Range("I3").Value = "Amul Butter (500g)"

$sales.Range("J3").NumberFormat = "@"
$sales.Range("J3").Value = "5000"

$sales.Range("K3").NumberFormat = "@"
$sales.Range("K3").Value = "4000"

$sales.Range("L3").Value = "2025-09-23T13:02:08.595Z"
